$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Bandage" (绷带) card's effect text is simplified: it no longer lets the
# player discard extra copies of itself for bonus healing, it now just heals 1.
$ws.Range("B9").Value = "回复1生命。"

# Update the active selection to match the author's final cursor position
# (B11) and drop the old scrolled viewport, matching the saved sheet view.
$ws.Range("B11").Select() | Out-Null
